$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "Lagt til prosjektfiler fra BB, og begynt å fylle inn systemdata - basecase"
# Rename the sheet to reflect the base-case load-flow scenario.
$ws.Name = "Base Case"

# Column A holds the line names ("1-2", "2-3", ...). Re-key them with the
# "n — m" (em dash) notation, store them as centred text.
$lineNames = @("1 — 2", "2 — 3", "3 — 4", "4 — 5", "5 — 6", "6 — 7", "7 — 8", "1 — 8", "1 — 6")

# Build the new text/centred format once on A2, then fan it out to A3:A10
# via copy/paste-special so every cell lands on the same style record.
$first = $ws.Cells.Item(2, 1)
$first.HorizontalAlignment = -4108   # xlCenter
$first.NumberFormat = "@"
$first.Copy()
for ($row = 3; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)   # xlPasteFormats
}

for ($i = 0; $i -lt $lineNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $lineNames[$i]
}

# Move the active selection (matches the author's last saved cursor position).
$ws.Range("C14").Select()
